$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.072131
$ws.Range("H2").Value = 18.216393
$ws.Range("I2").Value = 0.003943999267036455
$ws.Range("J2").Value = 0.003943999267036454
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 28.72417333333333
$ws.Range("N2").Value = 86.17251999999999
$ws.Range("O2").Value = 0.4233259107972328
$ws.Range("P2").Value = 0.4233259107972328
$ws.Range("Q2").Value = 174.4169433467066
$ws.Range("R2").Value = 1569.75249012036
$ws.Range("S2").Value = 0.001669597081901825
$ws.Range("T2").Value = 0.001669597081901825
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.072131
$ws.Range("H3").Value = 18.216393
$ws.Range("I3").Value = 0.003943999267036455
$ws.Range("J3").Value = 0.003943999267036454
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 30.56986233333333
$ws.Range("N3").Value = 91.709587
$ws.Range("O3").Value = 0.4505269713084062
$ws.Range("P3").Value = 0.4505269713084062
$ws.Range("Q3").Value = 185.6242087399657
$ws.Range("R3").Value = 1670.617878659691
$ws.Range("S3").Value = 0.001776878044620508
$ws.Range("T3").Value = 0.001776878044620508
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.072131
$ws.Range("H4").Value = 18.216393
$ws.Range("I4").Value = 0.003943999267036455
$ws.Range("J4").Value = 0.003943999267036454
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.559531999999999
$ws.Range("N4").Value = 25.678596
$ws.Range("O4").Value = 0.126147117894361
$ws.Range("P4").Value = 0.126147117894361
$ws.Range("Q4").Value = 51.97459960269199
$ws.Range("R4").Value = 467.771396424228
$ws.Range("S4").Value = 0.0004975241405141209
$ws.Range("T4").Value = 0.0004975241405141208
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1480.851806666667
$ws.Range("H5").Value = 4442.55542
$ws.Range("I5").Value = 0.9618498744646554
$ws.Range("J5").Value = 0.9618498744646552
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 28.72417333333333
$ws.Range("N5").Value = 86.17251999999999
$ws.Range("O5").Value = 0.4233259107972328
$ws.Range("P5").Value = 0.4233259107972328
$ws.Range("Q5").Value = 42536.24397567315
$ws.Range("R5").Value = 382826.1957810583
$ws.Range("S5").Value = 0.4071759741579543
$ws.Range("T5").Value = 0.4071759741579542
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1480.851806666667
$ws.Range("H6").Value = 4442.55542
$ws.Range("I6").Value = 0.9618498744646554
$ws.Range("J6").Value = 0.9618498744646552
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.56986233333333
$ws.Range("N6").Value = 91.709587
$ws.Range("O6").Value = 0.4505269713084062
$ws.Range("P6").Value = 0.4505269713084062
$ws.Range("Q6").Value = 45269.43586586795
$ws.Range("R6").Value = 407424.9227928115
$ws.Range("S6").Value = 0.433339310795932
$ws.Range("T6").Value = 0.4333393107959319
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1480.851806666667
$ws.Range("H7").Value = 4442.55542
$ws.Range("I7").Value = 0.9618498744646554
$ws.Range("J7").Value = 0.9618498744646552
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.559531999999999
$ws.Range("N7").Value = 25.678596
$ws.Range("O7").Value = 0.126147117894361
$ws.Range("P7").Value = 0.126147117894361
$ws.Range("Q7").Value = 12675.39842642114
$ws.Range("R7").Value = 114078.5858377903
$ws.Range("S7").Value = 0.1213345895107692
$ws.Range("T7").Value = 0.1213345895107692
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 52.663316
$ws.Range("H8").Value = 157.989948
$ws.Range("I8").Value = 0.03420612626830831
$ws.Range("J8").Value = 0.0342061262683083
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 28.72417333333333
$ws.Range("N8").Value = 86.17251999999999
$ws.Range("O8").Value = 0.4233259107972328
$ws.Range("P8").Value = 0.4233259107972328
$ws.Range("Q8").Value = 1512.710217092106
$ws.Range("R8").Value = 13614.39195382896
$ws.Range("S8").Value = 0.01448033955737676
$ws.Range("T8").Value = 0.01448033955737676
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 52.663316
$ws.Range("H9").Value = 157.989948
$ws.Range("I9").Value = 0.03420612626830831
$ws.Range("J9").Value = 0.0342061262683083
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.56986233333333
$ws.Range("N9").Value = 91.709587
$ws.Range("O9").Value = 0.4505269713084062
$ws.Range("P9").Value = 0.4505269713084062
$ws.Range("Q9").Value = 1609.910320136831
$ws.Range("R9").Value = 14489.19288123148
$ws.Range("S9").Value = 0.01541078246785386
$ws.Range("T9").Value = 0.01541078246785386
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 52.663316
$ws.Range("H10").Value = 157.989948
$ws.Range("I10").Value = 0.03420612626830831
$ws.Range("J10").Value = 0.0342061262683083
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.559531999999999
$ws.Range("N10").Value = 25.678596
$ws.Range("O10").Value = 0.126147117894361
$ws.Range("P10").Value = 0.126147117894361
$ws.Range("Q10").Value = 450.773338528112
$ws.Range("R10").Value = 4056.960046753008
$ws.Range("S10").Value = 0.004315004243077686
$ws.Range("T10").Value = 0.004315004243077685
